# Update the "Förändrad" (Changed) date column (C) for rows 2-51
# from serial date 45178 (2023-09-09) to 45179 (2023-09-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C51").Value = 45179
